$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2029.7142
$ws.Range("I19").Value = 1660
$ws.Range("K19").Value = 1660
$ws.Range("M19").Value = -1485

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3077
$ws.Range("I62").Value = 2624.25
$ws.Range("K62").Value = 2624.25
$ws.Range("M62").Value = -2000.25

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3077
$ws.Range("I65").Value = 2624.25
$ws.Range("K65").Value = 13121.25
$ws.Range("M65").Value = -10001.25

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3326.818
$ws.Range("J88").Value = 3261.111
$ws.Range("L88").Value = 3261.111
$ws.Range("N88").Value = -4073.111

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 3326.818
$ws.Range("J91").Value = 3261.111
$ws.Range("L91").Value = 3261.111
$ws.Range("N91").Value = -6069.111

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 10423.111
$ws.Range("I113").Value = 10543.286
$ws.Range("K113").Value = 10543.286
$ws.Range("M113").Value = -7289.286

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1061.7858
$ws.Range("I137").Value = 946.375
$ws.Range("J137").Value = 1215.6666
$ws.Range("K137").Value = 2839.125
$ws.Range("L137").Value = 3646.9998
$ws.Range("M137").Value = -289.125
$ws.Range("N137").Value = -8746.9998

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3413.0454
$ws.Range("I138").Value = 3036.625
$ws.Range("J138").Value = 3628.1428
$ws.Range("K138").Value = 9109.875
$ws.Range("L138").Value = 10884.4284
$ws.Range("M138").Value = -3969.875
$ws.Range("N138").Value = -21164.4284

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2783.7144
$ws.Range("I2").Value = 1581.3334
$ws.Range("K2").Value = 1581.3334
$ws.Range("M2").Value = -1468.3334

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3834.7144
$ws.Range("I45").Value = 1449.5
$ws.Range("J45").Value = 4788.8
$ws.Range("K45").Value = 1449.5
$ws.Range("L45").Value = 4788.8
$ws.Range("M45").Value = -1072.5
$ws.Range("N45").Value = -5542.8

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2783.7144
$ws.Range("I116").Value = 1581.3334
$ws.Range("K116").Value = 1581.3334
$ws.Range("M116").Value = 712.6666

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2783.7144
$ws.Range("I3").Value = 1581.3334
$ws.Range("K3").Value = 1581.3334
$ws.Range("M3").Value = -1467.3334

# BSM row 8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 724.5
$ws.Range("I8").Value = 724.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 724.5
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = -584.5

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2185.1667
$ws.Range("I20").Value = 1140.6666
$ws.Range("K20").Value = 1140.6666
$ws.Range("M20").Value = -893.6666

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1318
$ws.Range("I134").Value = 1403.6666
$ws.Range("K134").Value = 4210.9998
$ws.Range("M134").Value = -1675.9998

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 12136.889
$ws.Range("I105").Value = 13404
$ws.Range("K105").Value = 13404
$ws.Range("M105").Value = -11657

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 333.0476
$ws.Range("I107").Value = 302.35294
$ws.Range("J107").Value = 463.5
$ws.Range("K107").Value = 302.35294
$ws.Range("L107").Value = 463.5
$ws.Range("M107").Value = 1617.64706
$ws.Range("N107").Value = -4303.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3026.8125
$ws.Range("I132").Value = 2571.6924
$ws.Range("K132").Value = 7715.0772
$ws.Range("M132").Value = -5185.0772

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 769.2
$ws.Range("I5").Value = 711.75
$ws.Range("K5").Value = 2135.25
$ws.Range("M5").Value = -2023.25

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 567
$ws.Range("I26").Value = 812.5
$ws.Range("J26").Value = 76
$ws.Range("K26").Value = 2437.5
$ws.Range("L26").Value = 228
$ws.Range("M26").Value = -2149.5
$ws.Range("N26").Value = -804

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1502.5
$ws.Range("J68").Value = 1502.5
$ws.Range("L68").Value = 4507.5
$ws.Range("N68").Value = -6129.5

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1502.5
$ws.Range("J71").Value = 1502.5
$ws.Range("L71").Value = 13522.5
$ws.Range("N71").Value = -21634.5

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 282.42856
$ws.Range("I98").Value = 98
$ws.Range("J98").Value = 313.16666
$ws.Range("K98").Value = 294
$ws.Range("L98").Value = 939.4999799999999
$ws.Range("M98").Value = 1204
$ws.Range("N98").Value = -3935.49998

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 890.8125

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 769.2
$ws.Range("I135").Value = 711.75
$ws.Range("K135").Value = 6405.75
$ws.Range("M135").Value = -3870.75

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1200
$ws.Range("I140").Value = 1500
$ws.Range("J140").Value = 600
$ws.Range("K140").Value = 4500
$ws.Range("L140").Value = 1800
$ws.Range("M140").Value = 680
$ws.Range("N140").Value = -12160

# GSM row 40
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10302

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11667
$ws.Range("I70").Value = 11667
$ws.Range("K70").Value = 11667
$ws.Range("M70").Value = -11397

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11667
$ws.Range("I73").Value = 11667
$ws.Range("K73").Value = 11667
$ws.Range("M73").Value = -10731

# GSM row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 723.625
$ws.Range("I16").Value = 723.625
$ws.Range("K16").Value = 723.625
$ws.Range("M16").Value = -553.625

# LTW row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 19999
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# WVR row 32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 5265
$ws.Range("I32").Value = 5265
$ws.Range("K32").Value = 5265
$ws.Range("M32").Value = -4948

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3480.7
$ws.Range("I81").Value = 3423
$ws.Range("K81").Value = 6846
$ws.Range("M81").Value = -5785

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3480.7
$ws.Range("I84").Value = 3423
$ws.Range("K84").Value = 34230
$ws.Range("M84").Value = -28926

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2640.923
$ws.Range("I132").Value = 2343.2
$ws.Range("K132").Value = 7029.599999999999
$ws.Range("M132").Value = -4499.599999999999
